# Refresh the cryptocurrency price/volume snapshot (scrape update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.145.86"
$ws.Range("E2").Value = "  +4.01%  "
$ws.Range("D3").Value = "2.349.38"
$ws.Range("E3").Value = "  +2.41%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "521.69"
$ws.Range("E5").Value = "  +2.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.17"
$ws.Range("E6").Value = "  +4.25%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.48%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.538"
$ws.Range("E8").Value = "  +1.44%  "
$ws.Range("D9").Value = "2.348.86"
$ws.Range("E9").Value = "  +1.42%  "
$ws.Range("E10").Value = "  +5.86%  "
$ws.Range("E11").Value = "  -0.92%  "
$ws.Range("E12").Value = "  +4.08%  "
$ws.Range("E13").Value = "  +0.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.94"
$ws.Range("E14").Value = "  +0.47%  "
$ws.Range("D15").Value = "2.742.23"
$ws.Range("E15").Value = "  +1.51%  "
$ws.Range("D16").Value = "56.963.51"
$ws.Range("E16").Value = "  +3.68%  "
$ws.Range("E17").Value = "  +2.08%  "
$ws.Range("D18").Value = "2.329.03"
$ws.Range("E18").Value = "  +1.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.50"
$ws.Range("E19").Value = "  -0.18%  "
$ws.Range("E20").Value = "  +0.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "323.02"
$ws.Range("E21").Value = "  +3.73%  "
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.87"
$ws.Range("E24").Value = "  +0.96%  "
$ws.Range("E25").Value = "  +9.41%  "
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("E27").Value = "  +6.10%  "
$ws.Range("E28").Value = "  +13.61%  "
$ws.Range("D29").Value = "0.0₃0743"
$ws.Range("E29").Value = "  +4.61%  "
$ws.Range("E30").Value = "  +5.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "167.94"
$ws.Range("E31").Value = "  -2.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.22"
$ws.Range("E32").Value = "  +1.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.37"
$ws.Range("E33").Value = "  +1.65%  "
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.61%  "
$ws.Range("E36").Value = "  +1.71%  "
$ws.Range("E37").Value = "  +0.87%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.04"
$ws.Range("E38").Value = "  +4.01%  "
$ws.Range("E39").Value = "  +7.51%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.95"
$ws.Range("E40").Value = "  +3.05%  "
$ws.Range("E41").Value = "  +0.71%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.61"
$ws.Range("E42").Value = "  +4.81%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "138.30"
$ws.Range("E43").Value = "  +3.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.21"
$ws.Range("E44").Value = "  +5.73%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "277.53"
$ws.Range("E45").Value = "  +6.06%  "
$ws.Range("E46").Value = "  +2.22%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0505"
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("E48").Value = "  +2.38%  "
$ws.Range("E49").Value = "  +3.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.90"
$ws.Range("E50").Value = "  +8.13%  "
$ws.Range("E51").Value = "  +0.73%  "
